$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) column cells so numeric-looking strings
# like "315.50" or "1.90" are preserved exactly as text, not coerced to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Price (D) and Volume(1h) (E) column updates
$ws.Cells.Item(2, 4).Value = "45.409.06"
$ws.Cells.Item(2, 5).Value = "  -0.04%  "
$ws.Cells.Item(3, 4).Value = "2.375.75"
$ws.Cells.Item(3, 5).Value = "  -0.64%  "
$ws.Cells.Item(4, 5).Value = "  -0.08%  "
$ws.Cells.Item(5, 4).Value = "315.50"
$ws.Cells.Item(5, 5).Value = "  -0.93%  "
$ws.Cells.Item(6, 4).Value = "108.88"
$ws.Cells.Item(6, 5).Value = "  -3.85%  "
$ws.Cells.Item(7, 4).Value = "0.640"
$ws.Cells.Item(7, 5).Value = "  +0.30%  "
$ws.Cells.Item(8, 5).Value = "  +0.01%  "
$ws.Cells.Item(9, 4).Value = "0.616"
$ws.Cells.Item(9, 5).Value = "  -2.09%  "
$ws.Cells.Item(10, 4).Value = "40.98"
$ws.Cells.Item(10, 5).Value = "  -3.62%  "
$ws.Cells.Item(11, 5).Value = "  -1.64%  "
$ws.Cells.Item(12, 4).Value = "8.55"
$ws.Cells.Item(12, 5).Value = "  -1.88%  "
$ws.Cells.Item(13, 5).Value = "  +0.96%  "
$ws.Cells.Item(14, 4).Value = "0.985"
$ws.Cells.Item(14, 5).Value = "  -3.52%  "
$ws.Cells.Item(15, 4).Value = "2.737.12"
$ws.Cells.Item(15, 5).Value = "  -0.59%  "
$ws.Cells.Item(16, 4).Value = "15.45"
$ws.Cells.Item(16, 5).Value = "  -2.80%  "
$ws.Cells.Item(17, 4).Value = "2.368.33"
$ws.Cells.Item(17, 5).Value = "  -0.92%  "
$ws.Cells.Item(18, 4).Value = "45.387.22"
$ws.Cells.Item(18, 5).Value = "  -0.06%  "
$ws.Cells.Item(19, 4).Value = "15.89"
$ws.Cells.Item(19, 5).Value = "  +18.18%  "
$ws.Cells.Item(20, 4).Value = "7.34"
$ws.Cells.Item(20, 5).Value = "  -4.53%  "
$ws.Cells.Item(21, 5).Value = "  -1.69%  "
$ws.Cells.Item(22, 4).Value = "3.61"
$ws.Cells.Item(22, 5).Value = "  +1.62%  "
$ws.Cells.Item(23, 4).Value = "73.33"
$ws.Cells.Item(23, 5).Value = "  -2.34%  "
$ws.Cells.Item(24, 4).Value = "261.19"
$ws.Cells.Item(24, 5).Value = "  -3.26%  "
$ws.Cells.Item(25, 4).Value = "2.36"
$ws.Cells.Item(25, 5).Value = "  -0.68%  "
$ws.Cells.Item(26, 5).Value = "  +0.27%  "
$ws.Cells.Item(27, 4).Value = "7.64"
$ws.Cells.Item(27, 5).Value = "  +1.68%  "
$ws.Cells.Item(28, 5).Value = "  -1.10%  "
$ws.Cells.Item(29, 5).Value = "  -1.57%  "
$ws.Cells.Item(30, 4).Value = "0.0969"
$ws.Cells.Item(30, 5).Value = "  +2.67%  "
$ws.Cells.Item(31, 4).Value = "22.37"
$ws.Cells.Item(31, 5).Value = "  -2.49%  "
$ws.Cells.Item(32, 4).Value = "37.29"
$ws.Cells.Item(32, 5).Value = "  -4.56%  "
$ws.Cells.Item(33, 4).Value = "166.57"
$ws.Cells.Item(33, 5).Value = "  -2.32%  "
$ws.Cells.Item(34, 5).Value = "  -3.36%  "
$ws.Cells.Item(35, 4).Value = "0.133"
$ws.Cells.Item(35, 5).Value = "  -0.75%  "
$ws.Cells.Item(36, 4).Value = "0.117"
$ws.Cells.Item(36, 5).Value = "  -0.71%  "
$ws.Cells.Item(37, 4).Value = "4.73"
$ws.Cells.Item(37, 5).Value = "  -3.91%  "
$ws.Cells.Item(38, 4).Value = "4.05"
$ws.Cells.Item(38, 5).Value = "  +1.98%  "
$ws.Cells.Item(39, 4).Value = "1.90"
$ws.Cells.Item(39, 5).Value = "  +8.17%  "
$ws.Cells.Item(40, 4).Value = "2.93"
$ws.Cells.Item(40, 5).Value = "  -6.83%  "
$ws.Cells.Item(41, 4).Value = "0.0355"
$ws.Cells.Item(41, 5).Value = "  -2.83%  "
$ws.Cells.Item(42, 4).Value = "98.45"
$ws.Cells.Item(42, 5).Value = "  -6.17%  "
$ws.Cells.Item(43, 4).Value = "70.43"
$ws.Cells.Item(43, 5).Value = "  -2.05%  "
$ws.Cells.Item(44, 5).Value = "  -3.51%  "
$ws.Cells.Item(45, 4).Value = "0.228"
$ws.Cells.Item(45, 5).Value = "  -5.84%  "
$ws.Cells.Item(46, 5).Value = "  +0.02%  "
$ws.Cells.Item(47, 4).Value = "5.97"
$ws.Cells.Item(47, 5).Value = "  +2.68%  "
$ws.Cells.Item(48, 4).Value = "1.819.79"
$ws.Cells.Item(48, 5).Value = "  +9.86%  "
$ws.Cells.Item(49, 4).Value = "83.97"
$ws.Cells.Item(49, 5).Value = "  +6.90%  "
$ws.Cells.Item(50, 4).Value = "111.65"
$ws.Cells.Item(50, 5).Value = "  -5.27%  "
$ws.Cells.Item(51, 4).Value = "9.25"
$ws.Cells.Item(51, 5).Value = "  -1.47%  "
